$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.896.65"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$ws.Range("D3").Value = "1.814.11"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.97"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4730"
$ws.Range("E7").Value = "  +2.50%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3701"
$ws.Range("E8").Value = "  -1.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07377"
$ws.Range("E9").Value = "  -0.46%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8694"
$ws.Range("E10").Value = "  +0.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.40"
$ws.Range("E11").Value = "  -1.05%  "

# Row 12
$ws.Range("D12").Value = "1.867.46"
$ws.Range("E12").Value = "  +2.67%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.347"
$ws.Range("E13").Value = "  -1.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07065"

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.79"
$ws.Range("E15").Value = "  -0.36%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.502"
$ws.Range("E16").Value = "  -2.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008713"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.69"
$ws.Range("E20").Value = "  -1.57%  "

# Row 21
$ws.Range("D21").Value = "26.874.15"
$ws.Range("E21").Value = "  -1.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.332"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  -3.35%  "

# Row 24
$ws.Range("D24").Value = "2.111.42"
$ws.Range("E24").Value = "  +3.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.897"
$ws.Range("E25").Value = "  -2.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.71"
$ws.Range("E26").Value = "  -0.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.125"
$ws.Range("E28").Value = "  -5.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.284"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.36"
$ws.Range("E30").Value = "  -1.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08958"
$ws.Range("E31").Value = "  +0.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7549"
$ws.Range("E32").Value = "  -3.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.153"
$ws.Range("E33").Value = "  -2.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.922"
$ws.Range("E34").Value = "  +0.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.451"
$ws.Range("E35").Value = "  -1.73%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.091"
$ws.Range("E37").Value = "  -1.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01950"
$ws.Range("E38").Value = "  -0.58%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05256"
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.923"
$ws.Range("E40").Value = "  +0.39%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5329"
$ws.Range("E41").Value = "  +0.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.190"
$ws.Range("E42").Value = "  -1.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.360"
$ws.Range("E43").Value = "  -1.43%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1654"
$ws.Range("E44").Value = "  -1.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.434"
$ws.Range("E45").Value = "  -1.95%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4941"
$ws.Range("E46").Value = "  -2.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("E47").Value = "  -1.05%  "

# Row 48
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.669"
$ws.Range("E49").Value = "  -0.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.14"
$ws.Range("E50").Value = "  -1.93%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06279"
$ws.Range("E51").Value = "  -0.75%  "
